$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 (RM 232) - everything below shifts up
$ws.Rows.Item(26).Delete()

# After the first deletion, the row that was "SC 92" (originally row 28)
# is now at row 27. Delete it too.
$ws.Rows.Item(27).Delete()
